$wb = $excel.ActiveWorkbook

# --- 1. Rename "Non-Hispanic Voters" sheet to the (31-char-limited) new name ---
$nonHispSheet = $wb.Worksheets.Item("Non-Hispanic Voters")
$nonHispSheet.Name = "Non-Hispanic and Non-Undesign"

# --- 2. Update header + values on that sheet to reflect "Non-Hispanic and Non-Undesignated Voters" ---
$nonHispSheet.Range("B1").Value = "Number of Non-Hispanic and Non-Undesignated Voters"

$nonHispSheet.Range("B2").Value = 2011705
$nonHispSheet.Range("B3").Value = 1387
$nonHispSheet.Range("B4").Value = 34603
$nonHispSheet.Range("B5").Value = 8240
$nonHispSheet.Range("B6").Value = 1914504
$nonHispSheet.Range("B7").Value = 2004857

# --- 3. Add a new "Undesignated Voters" sheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$newSheet.Name = "Undesignated Voters"

$newSheet.Range("A1").Value = "Party"
$newSheet.Range("B1").Value = "Number of Undesignated Voters"

$newSheet.Range("A2").Value = "DEM"
$newSheet.Range("B2").Value = 725368

$newSheet.Range("A3").Value = "GRE"
$newSheet.Range("B3").Value = 977

$newSheet.Range("A4").Value = "LIB"
$newSheet.Range("B4").Value = 17525

$newSheet.Range("A5").Value = "NLB"
$newSheet.Range("B5").Value = 5459

$newSheet.Range("A6").Value = "REP"
$newSheet.Range("B6").Value = 625600

$newSheet.Range("A7").Value = "UNA"
$newSheet.Range("B7").Value = 1014761

# Match formatting of header row on the other sheets (bold, centered - style index 1)
# by copying the existing header format rather than re-deriving it property by
# property (keeps the style table identical to the sibling sheets).
$lastSheet.Range("A1:B1").Copy()
$newSheet.Range("A1:B1").PasteSpecial(-4122) # xlPasteFormats
